# FIX: start level list in -1
#
# The numbered list (numId=1) that contains "primeiro item" ... "teste 2"
# was accidentally created one level too deep. Shift every paragraph of
# that sub-list up by one list level (ilvl -= 1 in OOXML terms, which is
# ListLevelNumber -= 1 in the 1-based Word object model).

$d = $word.ActiveDocument

$targets = @(
    "primeiro item",
    "segundo item",
    "terceiro item",
    "quarto item",
    "teste",
    "teste 1",
    "teste 2"
)

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text.Trim()
    if ($targets -contains $text) {
        $current = $p.Range.ListFormat.ListLevelNumber
        $p.Range.ListFormat.ListLevelNumber = $current - 1
    }
}
